$d = $word.ActiveDocument

# --- Paragraph 1 ("**ID__AFFARS_5319_topic_2__ID**" banner line) ---
$p1 = $d.Paragraphs(1)

# Add a thin paragraph border (top/left/bottom/right, 5pt space to text)
$p1.Range.Borders.DistanceFromTop = 5
$p1.Range.Borders.DistanceFromLeft = 5
$p1.Range.Borders.DistanceFromBottom = 5
$p1.Range.Borders.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p1.Format.LeftIndent = 11.25

# The paragraph currently holds two runs: the placeholder token, then a
# lone trailing space run. Drop the trailing space run entirely, then
# update the placeholder token text in place (same run, same rPr).
$pStart = $p1.Range.Start
$oldToken = "**ID__AFFARS_5319_topic_2__ID**"
$newToken = "**ID__AFFARS_SUBPART_5319_2__ID**"

$spaceRun = $d.Range($pStart + $oldToken.Length, $pStart + $oldToken.Length + 1)
$spaceRun.Delete()

$tokenRun = $d.Range($pStart, $pStart + $oldToken.Length)
$tokenRun.Text = $newToken

Write-Output $d.Content.Text
